$d = $word.ActiveDocument

# Locate the "PERSPECTIVE" heading paragraph that the new Exposure & Marketing
# Time section gets appended after.
$findRng = $d.Content
$found = $findRng.Find.Execute("PERSPECTIVE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'PERSPECTIVE' paragraph to anchor the insertion."
}

# Insert point: end of the document content (the "PERSPECTIVE" paragraph is
# the last paragraph in the body, right before the sectPr).
$insertRng = $d.Content
$insertRng.Collapse(0)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblW w:w="9715" w:type="dxa"/>
    <w:tblInd w:w="5" w:type="dxa"/>
    <w:shd w:val="pct20" w:color="auto" w:fill="auto"/>
    <w:tblLook w:val="01E0" w:firstRow="1" w:lastRow="1" w:firstColumn="1" w:lastColumn="1" w:noHBand="0" w:noVBand="0"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="2245"/>
    <w:gridCol w:w="7470"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="2245" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
        <w:hideMark/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="center" w:pos="4680"/>
          </w:tabs>
          <w:spacing w:before="120" w:after="120"/>
          <w:rPr>
            <w:b/>
            <w:bCs/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:bCs/>
          </w:rPr>
          <w:lastRenderedPageBreak/>
          <w:t>EXPOSURE AND MARKETING TIME</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="7470" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="auto"/>
        <w:hideMark/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="left" w:pos="4500"/>
          </w:tabs>
          <w:spacing w:before="120" w:after="120" w:line="280" w:lineRule="atLeast"/>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">Exposure </w:t>
        </w:r>
        <w:proofErr w:type="gramStart"/>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>time as</w:t>
        </w:r>
        <w:proofErr w:type="gramEnd"/>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve"> used in this appraisal report is defined as:</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="left" w:pos="4500"/>
          </w:tabs>
          <w:spacing w:before="120" w:after="120" w:line="280" w:lineRule="atLeast"/>
          <w:rPr>
            <w:i/>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>Exposure time.</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:i/>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve"> The estimated length of time the property interest being appraised would have been offered on the market prior to the hypothetical consummation of a sale at market value on the effective date of the appraisal; a retrospective estimate based on an analysis of past events assuming a competitive and open market.</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="left" w:pos="4500"/>
          </w:tabs>
          <w:spacing w:before="120" w:after="120" w:line="280" w:lineRule="atLeast"/>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">Based upon the market data from the sales of comparable properties in the </w:t>
        </w:r>
        <w:proofErr w:type="gramStart"/>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>market</w:t>
        </w:r>
        <w:proofErr w:type="gramEnd"/>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve"> my conclusion of exposure time follows:</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="left" w:pos="4500"/>
          </w:tabs>
          <w:spacing w:before="120" w:after="120" w:line="280" w:lineRule="atLeast"/>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>Conclusion:</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="left" w:pos="3672"/>
          </w:tabs>
          <w:spacing w:before="120" w:after="120" w:line="280" w:lineRule="atLeast"/>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>Exposure Time</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:tab/>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>[ENTER TIME HERE]</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="left" w:pos="3672"/>
          </w:tabs>
          <w:spacing w:before="120" w:after="120" w:line="280" w:lineRule="atLeast"/>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>Marketing time.</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:i/>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve"> The time it takes an interest in real property to sell on the market </w:t>
        </w:r>
        <w:proofErr w:type="gramStart"/>
        <w:r>
          <w:rPr>
            <w:i/>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>subsequent to</w:t>
        </w:r>
        <w:proofErr w:type="gramEnd"/>
        <w:r>
          <w:rPr>
            <w:i/>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve"> the date of appraisal</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:i/>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>.</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="left" w:pos="4500"/>
          </w:tabs>
          <w:spacing w:before="120" w:after="120" w:line="280" w:lineRule="atLeast"/>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve">Based upon the market data from the sales of comparable properties in the </w:t>
        </w:r>
        <w:proofErr w:type="gramStart"/>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>market</w:t>
        </w:r>
        <w:proofErr w:type="gramEnd"/>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t xml:space="preserve"> my conclusion of exposure time follows:</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="left" w:pos="4500"/>
          </w:tabs>
          <w:spacing w:before="120"/>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>Conclusion:</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="left" w:pos="3672"/>
          </w:tabs>
          <w:spacing w:before="120" w:after="120" w:line="280" w:lineRule="atLeast"/>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>Marketing Time</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:tab/>
        </w:r>
        <w:r>
          <w:rPr>
            <w:lang w:bidi="en-US"/>
          </w:rPr>
          <w:t>[ENTER TIME HERE]</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:tabs>
            <w:tab w:val="left" w:pos="3672"/>
          </w:tabs>
          <w:spacing w:before="120" w:after="120" w:line="280" w:lineRule="atLeast"/>
        </w:pPr>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
</w:p>
'@

[void]$insertRng.InsertXML($xml)
Write-Output "Inserted Exposure and Marketing Time section (paragraph + table + paragraph) after 'PERSPECTIVE'."
